# Add a parallel "Invalido" block (columns H:J) mirroring the existing
# Teste / Caracter Informado / Validação layout used for "Valido" in E:G,
# and a new "Mensagem" column (K) that reports a specific error message
# for each invalid test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers H1:J1 (copy format from E1:G1, same header style) ---
$ws.Range("E1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Teste"
$ws.Range("I1").Value = "Caracter Informado"
$ws.Range("J1").Value = "Validação"

# --- Row 2 (H2:J2) - moved from the old E5:G5 "Invalido" row ---
$ws.Range("E2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("H2").Value = "ABCDEFGHIJKMNOPQRSUVX"
$ws.Range("I2").Value = "A"
$ws.Range("J2").Value = "Invalido"

# --- Row 3 (H3:J3) - moved from the old E6:G6 "Invalido" row ---
$ws.Range("E3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("H3").Value = "ABCDEFGHIJKMNOPQ"
$ws.Range("I3").Value = "U"
$ws.Range("J3").Value = "Invalido"

# --- Row 4 (H4:J4) - moved from the old E7:G7 "Invalido" row ---
$ws.Range("E4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("H4").Value = "ABCDEFGHIJKMNOPQ"
$ws.Range("I4").Value = "Não informado"
$ws.Range("J4").Value = "Invalido"

# --- New "Mensagem" column K ---
# K1: header - same dark/white header look as A1/E1, but only left/right
# thin borders (no top/bottom) so it blends with the red message cells below.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Mensagem"
$ws.Range("K1").Borders(8).LineStyle = -4142
$ws.Range("K1").Borders(9).LineStyle = -4142
$ws.Range("K1").Borders(7).LineStyle = 1
$ws.Range("K1").Borders(10).LineStyle = 1

$redColor = $ws.Range("G5").Interior.Color

# K2: error message for the "more than 20 chars" case (red fill, left/right border)
$ws.Range("K2").Value = "Erro Mais de 20"
$ws.Range("K2").Interior.Color = $redColor
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").Borders(7).LineStyle = 1
$ws.Range("K2").Borders(10).LineStyle = 1

# K3: no message configured for this case - just the red fill, no border/alignment
$ws.Range("K3").Value = -1
$ws.Range("K3").Interior.Color = $redColor

# K4: error message for the "not informed" case (red fill, left/right border)
$ws.Range("K4").Value = "Erro Não informado"
$ws.Range("K4").Interior.Color = $redColor
$ws.Range("K4").HorizontalAlignment = -4108
$ws.Range("K4").Borders(7).LineStyle = 1
$ws.Range("K4").Borders(10).LineStyle = 1

# --- Remove the old E5:G7 "Invalido" block, now represented in H:K ---
$ws.Range("E5:G5").Clear()
$ws.Rows("6:7").Delete()

# --- Column widths: best-effort resize of the touched columns ---
$ws.Columns("A").ColumnWidth = 25.140625
$ws.Columns("B").ColumnWidth = 10
$ws.Columns("C").ColumnWidth = 9.140625
$ws.Columns("E").ColumnWidth = 26.42578125
$ws.Columns("F").ColumnWidth = 18.28515625
$ws.Columns("H").ColumnWidth = 26.42578125
$ws.Columns("I").ColumnWidth = 17.28515625
$ws.Columns("J").ColumnWidth = 10.28515625
$ws.Columns("K").ColumnWidth = 19.140625

# --- Selection, matching the saved view state ---
$ws.Range("K2:K4").Select()

Write-Host "edit applied"
